$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuilt "Estado de Cuenta" detail table (rows 16-57): the underlying
# macro/database was regenerated, so each worker's block of periods is now
# grouped together (most recent period 2207 first, descending), while the
# set of (worker, periodo, valor, salario) tuples is unchanged.
$rows = @(
    @{ Row=16; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2207'; Valor=30284; Salario=908526 }
    @{ Row=17; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2206'; Valor=36341; Salario=908526 }
    @{ Row=18; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2205'; Valor=36341; Salario=908526 }
    @{ Row=19; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2204'; Valor=36341; Salario=908526 }
    @{ Row=20; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2203'; Valor=36341; Salario=908526 }
    @{ Row=21; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2202'; Valor=36341; Salario=908526 }
    @{ Row=22; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2201'; Valor=36341; Salario=908526 }
    @{ Row=23; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2112'; Valor=36341; Salario=908526 }
    @{ Row=24; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2111'; Valor=36341; Salario=908526 }
    @{ Row=25; Doc='79568200'; Nombre='JORGE ENRIQUE VISBAL MORENO'; Periodo='2110'; Valor=36341; Salario=908526 }
    @{ Row=26; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2207'; Valor=140000; Salario=2100000 }
    @{ Row=27; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2206'; Valor=84000; Salario=2100000 }
    @{ Row=28; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2205'; Valor=84000; Salario=2100000 }
    @{ Row=29; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2204'; Valor=84000; Salario=2100000 }
    @{ Row=30; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2203'; Valor=84000; Salario=2100000 }
    @{ Row=31; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2202'; Valor=168000; Salario=2100000 }
    @{ Row=32; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2201'; Valor=168000; Salario=2100000 }
    @{ Row=33; Doc='1047455106'; Nombre='DAVID PAREDES SEGOVIA'; Periodo='2112'; Valor=168000; Salario=2100000 }
    @{ Row=34; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2207'; Valor=30284; Salario=908526 }
    @{ Row=35; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2206'; Valor=36341; Salario=908526 }
    @{ Row=36; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2205'; Valor=36341; Salario=908526 }
    @{ Row=37; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2204'; Valor=36341; Salario=908526 }
    @{ Row=38; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2203'; Valor=36341; Salario=908526 }
    @{ Row=39; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2202'; Valor=36341; Salario=908526 }
    @{ Row=40; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2201'; Valor=36341; Salario=908526 }
    @{ Row=41; Doc='1047474476'; Nombre='JUAN CARLOS SARMIENTO ROSSO'; Periodo='2112'; Valor=36341; Salario=908526 }
    @{ Row=42; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2207'; Valor=30284; Salario=908526 }
    @{ Row=43; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2206'; Valor=36341; Salario=908526 }
    @{ Row=44; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2205'; Valor=36341; Salario=908526 }
    @{ Row=45; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2204'; Valor=36341; Salario=908526 }
    @{ Row=46; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2203'; Valor=36341; Salario=908526 }
    @{ Row=47; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2202'; Valor=36341; Salario=908526 }
    @{ Row=48; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2201'; Valor=36341; Salario=908526 }
    @{ Row=49; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2112'; Valor=36341; Salario=908526 }
    @{ Row=50; Doc='1047470413'; Nombre='ANGELICA MARIA LORDUY JIMENEZ'; Periodo='2109'; Valor=36341; Salario=908526 }
    @{ Row=51; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2207'; Valor=33333; Salario=3000000 }
    @{ Row=52; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2206'; Valor=40000; Salario=3000000 }
    @{ Row=53; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2205'; Valor=40000; Salario=3000000 }
    @{ Row=54; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2204'; Valor=40000; Salario=3000000 }
    @{ Row=55; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2203'; Valor=40000; Salario=3000000 }
    @{ Row=56; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2202'; Valor=40000; Salario=3000000 }
    @{ Row=57; Doc='1002308980'; Nombre='ELEANA MAGDALENA HURTADO BABILONIA'; Periodo='2201'; Valor=40000; Salario=3000000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc        # C: N Doc Trabajador
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre     # D: Nombre Trabajador
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo    # E: Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor      # F: Valor Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario    # G: Salario Basico
}
